$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# --- Update the time_taken values on the "data" sheet (column F, rows 2-7) ---
$dataSheet.Range("F2").Value = "2021-10-05 14:20:07.426701"
$dataSheet.Range("F3").Value = "2021-10-05 14:20:07.426709"
$dataSheet.Range("F4").Value = "2021-10-05 14:20:07.426712"
$dataSheet.Range("F5").Value = "2021-10-05 14:20:07.426715"
$dataSheet.Range("F6").Value = "2021-10-05 14:20:07.426718"
$dataSheet.Range("F7").Value = "2021-10-05 14:20:07.426720"

# --- Add the new "metadata" sheet right after "data" ---
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$newSheet.Name = "metadata"

# Header row (B1:G1) - same bold/bordered style as the "data" sheet header
$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

$dataSheet.Range("B1:F1").Copy()
$newSheet.Range("B1:F1").PasteSpecial(-4122)
$dataSheet.Range("B1").Copy()
$newSheet.Range("G1").PasteSpecial(-4122)

# Data row (A2:G2)
$newSheet.Range("A2").Value = 0
$dataSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

$newSheet.Range("B2").Value = "Endometrial cancer pertinent cancer susceptibility"
$newSheet.Range("C2").Value = 271
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "1.0"
# Clear the text-format style back off D2 so it matches the unstyled source cell
$dataSheet.Range("Z1").Copy()
$newSheet.Range("D2").PasteSpecial(-4122)
$newSheet.Range("E2").Value = "2017-11-05T02:37:20.335660Z"
$newSheet.Range("F2").Value = "2021-10-05 14:20:07.423411"
$newSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/271/?format=json"

$excel.CutCopyMode = $false

# Keep "data" as the active sheet/tab (matches the original activeTab="0")
$dataSheet.Activate()
